# Apply the "table column classification" edit to the fond.xlsx template.
#
# Summary of the change:
#   - Row 7 holds the table's placeholder header tokens used by the
#     templating engine. Three of them are being renamed/reclassified
#     from positional references ({child[0]:linked}, {child[1]:linked},
#     {child[2]}) to named-column references ({child[ID]:linked},
#     {child[DESCRIPTION]:linked}, {child[DATE]}).
#   - The active selection moves from P2:P3 to the single cell D7, and
#     the sheet view no longer needs to be scrolled (topLeftCell reset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reclassify the three positional child-column placeholders to named ones.
$ws.Range("A7").Value = "{child[ID]:linked}"
$ws.Range("B7").Value = "{child[DESCRIPTION]:linked}"
$ws.Range("C7").Value = "{child[DATE]}"

# Update the active selection/scroll position shown in the sheet view.
$ws.Range("D7").Select()
